$d = $word.ActiveDocument

# The document ends with: "... Kale 羽衣甘藍" paragraph, followed by a
# trailing (otherwise empty) paragraph that only carries the _GoBack
# bookmark. We need to add three new word/translation paragraphs right
# after "Kale" and before that trailing bookmark paragraph:
#   Chinese cabbage 白菜
#   Chinese flowering cabbage 菜心
#   Give in 屈服
#
# Replace the last paragraph's range with: the three new paragraphs
# followed by the original trailing bookmark paragraph (so the
# bookmark itself is preserved unchanged).

$paras = $d.Paragraphs
$n = $paras.Count
$lastPara = $paras.Item($n)
$r = $lastPara.Range

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>C</w:t></w:r><w:r><w:t xml:space="preserve">hinese cabbage </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>白菜</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>C</w:t></w:r><w:r><w:t xml:space="preserve">hinese flowering cabbage </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>菜心</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>G</w:t></w:r><w:r><w:t xml:space="preserve">ive in </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>屈服</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="001A0967" w:rsidRPr="003E7C52" w:rsidRDefault="001A0967"><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$r.InsertXML($xml)

Write-Output ("Paragraphs now: " + $d.Paragraphs.Count)
